# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Values in columns D (Price) and E (Volume(1h)) are authored as literal text
# in this sheet, so numeric-looking replacements are entered with a leading
# apostrophe to stop Excel's COM layer from auto-coercing them into real
# numbers (which would also silently drop significant trailing zeros, e.g.
# "8.90" -> 8.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.947.90'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '2.498.09'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''591.38'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '''175.03'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '2.497.59'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +7.10%  '
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = '''4.97'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = '2.956.67'
$ws.Range('D15').Value = '''25.56'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '68.820.94'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = '2.500.70'
$ws.Range('E18').Value = '  -3.26%  '
$ws.Range('D19').Value = '''359.23'
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').Value = '''7.51'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').Value = '''69.81'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('D25').Value = '''4.16'
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('D26').Value = '''8.90'
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('E27').Value = '  -6.46%  '
$ws.Range('D28').Value = '2.626.26'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('D30').Value = '''507.51'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '0.0₃0876'
$ws.Range('E31').Value = '  -3.05%  '
$ws.Range('D32').Value = '''7.69'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '''1.76'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '''1.20'
$ws.Range('E34').Value = '  -4.33%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '''162.48'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  -3.73%  '
$ws.Range('D38').Value = '''18.56'
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('D39').Value = '''18.67'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('E42').Value = '  -3.08%  '
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('E44').Value = '  -3.78%  '
$ws.Range('E45').Value = '  -5.33%  '
$ws.Range('D46').Value = '''148.88'
$ws.Range('E46').Value = '  +2.61%  '
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('D48').Value = '''0.508'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('E50').Value = '  -2.45%  '
$ws.Range('D51').Value = '''0.574'
$ws.Range('E51').Value = '  -2.21%  '
